$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet
$ws.Name = "Sheet1"

# Update the tab-label cell and remove the now-unused trailing row
$ws.Range("A2").Value = "Programstab"
$ws.Rows.Item(4).Delete()

# Resize columns to the new layout
$ws.Columns.Item(1).ColumnWidth = 13
$ws.Columns.Item(2).ColumnWidth = 68
$ws.Columns.Item(3).ColumnWidth = 60.166666666666664
$ws.Columns.Item(4).ColumnWidth = 23.833333333333332
$ws.Columns.Item(5).ColumnWidth = 29.166666666666668

# Resize the header / data rows
$ws.Rows.Item(1).RowHeight = 30.75
$ws.Rows.Item(2).RowHeight = 266.25

# Move the active selection
$ws.Range("C2").Select() | Out-Null
